$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 16 ("Unión Q-Learning y regrésión lineal"), pushing
#    the old rows 16-22 (Servicios en línea ... Implementación con JPA e Hibernate)
#    down to rows 17-23. Row 15's formatting (style 3 on col A, style 1 on the
#    hour columns) is inherited automatically by Insert().
$ws.Rows(16).Insert()

# New task row under "Implementación aprendizaje automático"
$ws.Range("A16").Value = "Unión Q-Learning y regrésión lineal"
$ws.Range("AK16").Value = "2.5 h."

# 2. Extend the date header (row 2) with two more days, copying the date
#    cell style (s=4) from the last existing date column (AJ2).
$ws.Range("AJ2").Copy() | Out-Null
$ws.Range("AK2:AL2").PasteSpecial(-4122) | Out-Null
$ws.Range("AK2").Value = 44038
$ws.Range("AL2").Value = 44039

# 3. New hour entries on existing rows, copying the "hours cell" style (s=1)
#    from a neighboring cell on the same row before writing the value.
$ws.Range("AH5").Copy() | Out-Null
$ws.Range("AL5").PasteSpecial(-4122) | Out-Null
$ws.Range("AL5").Value = "3 h."

$ws.Range("AJ15").Copy() | Out-Null
$ws.Range("AK15:AL15").PasteSpecial(-4122) | Out-Null
$ws.Range("AK15").Value = "2.5 h."
$ws.Range("AL15").Value = "1.5 h"

# 4. Update the totals text.
$ws.Range("A1").Value = "Total horas: 146.5"
$ws.Range("A25").Value = "4+3+1.5+4+4+3+2+3.5+4.5+4.5+5+1+5+3+4+5+4+4+3+4+3.5+8+3.5+5+5+5+4+3+3.5+4+2.5+5.5+6+6.5+5+4.5"

# 5. Restore the active selection.
$ws.Range("K20").Select()

$excel.CutCopyMode = 0
